$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2/C2 were text "-" placeholders -> now literal 0 numbers
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# New headers: Area (G1), Atotal (H1)
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# New Area column formulas
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Atotal sum
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Summary block in J/K
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Match the new selection from the diff
$ws.Range("J2:K2").Select()

$wb.Application.Calculate()
